{"js": "// Replace every occurrence of the literal \"{{else}}\" with \"{{#else}}\"\n// throughout the document body (covers both the stand-alone {{else}}\n// paragraphs and the inline {{#if GiftWrap}}Yes{{else}}No{{/if}} usage).\nconst body = context.document.body;\nconst results = body.search(\"{{else}}\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{{#else}}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace every occurrence of the literal \"{{else}}\" with \"{{#else}}\"\n# throughout the document (covers both the stand-alone {{else}} paragraphs\n# and the inline {{#if GiftWrap}}Yes{{else}}No{{/if}} usage).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{{else}}\"\n$find.Replacement.Text = \"{{#else}}\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n\n# wdReplaceAll = 2\n$find.Execute([ref]\"{{else}}\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"{{#else}}\", [ref]2) | Out-Null\n"}
